$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.357.73'
$ws.Range("E2").Value = '  +2.65%  '
$ws.Range("D3").Value = '2.423.29'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.34%  '
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +6.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.56'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0803'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("E12").Value = '  -2.05%  '
$ws.Range("E13").Value = '  -2.67%  '
$ws.Range("E14").Value = '  +1.81%  '
$ws.Range("D15").Value = '2.802.88'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").Value = '2.425.29'
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '45.259.02'
$ws.Range("E18").Value = '  +2.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.22%  '
$ws.Range("E20").Value = '  -1.64%  '
$ws.Range("D21").Value = '0.0₃0923'
$ws.Range("E21").Value = '  +1.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.86'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.66%  '
$ws.Range("E24").Value = '  -1.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.89%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.03%  '
$ws.Range("E29").Value = '  -11.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '49.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.124'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.92%  '
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0766'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("E37").Value = '  -2.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.45'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("E39").Value = '  -1.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.57%  '
$ws.Range("E41").Value = '  -3.32%  '
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.04%  '
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("D45").Value = '1.929.99'
$ws.Range("E45").Value = '  -1.43%  '
$ws.Range("E47").Value = '  +1.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.31%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +9.12%  '
$ws.Range("E50").Value = '  +4.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.56%  '
